# Add data for 2021-12-03: extends the "through November 24" snapshot to
# "through November 25" and bumps the November 2021 counts (column B) for
# the neighborhoods that saw a new carjacking recorded, plus some updates
# scattered in other month columns that were corrected at the same time.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet and update the title text in B1-equivalent (shared string 0)
$ws.Name = "Through 2021-11-25"
$ws.Range("B1").Value = "November 2021 (through November 25)"

# --- Cell updates (existing cells whose value changed) ---
$ws.Range("M2").Value = 17    # North Lawndale / November 2020
$ws.Range("X2").Value = 5     # North Lawndale / November 2019

$ws.Range("X3").Value = 5     # Garfield Park / November 2019
$ws.Range("AT3").Value = 6    # Garfield Park / November 2017

$ws.Range("BE7").Value = 4    # South Shore / November 2016

$ws.Range("AT8").Value = 5    # Englewood / November 2017
$ws.Range("BE8").Value = 5    # Englewood / November 2016

$ws.Range("B12").Value = 4    # Grand Boulevard / November 2021
$ws.Range("M12").Value = 7    # Grand Boulevard / November 2020

$ws.Range("M15").Value = 5    # Calumet Heights / November 2020

$ws.Range("B27").Value = 2    # South Chicago / November 2021

$ws.Range("B29").Value = 2    # Near South Side / November 2021

$ws.Range("B35").Value = 3    # Edgewater / November 2021

$ws.Range("B68").Value = 5    # Douglas / November 2021

$ws.Range("M96").Value = 3    # Ukrainian Village / November 2020

# --- New cell values (previously empty cells) ---
$ws.Range("AT16").Value = 1   # Washington Heights / November 2017
$ws.Range("M45").Value = 1    # United Center / November 2020
$ws.Range("M56").Value = 1    # Mount Greenwood / November 2020
$ws.Range("AI57").Value = 1   # Magnificent Mile / November 2018
$ws.Range("AI65").Value = 1   # Brighton Park / November 2018
$ws.Range("B99").Value = 1    # Wrigleyville / November 2021
